$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D and E on the data rows hold text that looks numeric/percent
# (e.g. "69.645.95", "3.30", "  +2.37%  "). Force text format first so
# Excel does not silently coerce these into numbers, then restore the
# original (default) style afterwards so no stray formatting is left behind.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '69.645.95'
$ws.Range('E2').Value = '  +2.37%  '
$ws.Range('D3').Value = '3.400.61'
$ws.Range('E3').Value = '  +1.98%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = '586.63'
$ws.Range('E5').Value = '  +0.49%  '
$ws.Range('D6').Value = '180.97'
$ws.Range('E6').Value = '  +2.97%  '
$ws.Range('D7').Value = '0.599'
$ws.Range('E7').Value = '  +0.99%  '
$ws.Range('D9').Value = '0.203'
$ws.Range('E9').Value = '  +11.34%  '
$ws.Range('D10').Value = '0.593'
$ws.Range('E10').Value = '  +1.93%  '
$ws.Range('D11').Value = '48.52'
$ws.Range('E11').Value = '  +1.85%  '
$ws.Range('D12').Value = '0.0000287'
$ws.Range('E12').Value = '  +5.06%  '
$ws.Range('D13').Value = '684.87'
$ws.Range('E13').Value = '  -2.07%  '
$ws.Range('D14').Value = '8.68'
$ws.Range('E14').Value = '  +3.37%  '
$ws.Range('D15').Value = '3.961.94'
$ws.Range('D16').Value = '69.723.95'
$ws.Range('E16').Value = '  +2.44%  '
$ws.Range('D17').Value = '3.414.37'
$ws.Range('E17').Value = '  +2.32%  '
$ws.Range('E18').Value = '  +1.51%  '
$ws.Range('D19').Value = '17.77'
$ws.Range('E19').Value = '  +1.70%  '
$ws.Range('D20').Value = '11.35'
$ws.Range('E20').Value = '  +1.89%  '
$ws.Range('D21').Value = '0.913'
$ws.Range('E21').Value = '  +1.84%  '
$ws.Range('D22').Value = '17.32'
$ws.Range('E22').Value = '  +1.94%  '
$ws.Range('D23').Value = '5.36'
$ws.Range('E23').Value = '  -0.38%  '
$ws.Range('D24').Value = '103.03'
$ws.Range('E24').Value = '  +2.10%  '
$ws.Range('D25').Value = '3.95'
$ws.Range('E25').Value = '  +0.82%  '
$ws.Range('D26').Value = '2.72'
$ws.Range('E26').Value = '  +1.33%  '
$ws.Range('D27').Value = '9.76'
$ws.Range('E27').Value = '  +3.70%  '
$ws.Range('D28').Value = '34.02'
$ws.Range('E28').Value = '  +2.40%  '
$ws.Range('D29').Value = '8.84'
$ws.Range('E29').Value = '  +3.53%  '
$ws.Range('D30').Value = '6.97'
$ws.Range('E30').Value = '  -0.39%  '
$ws.Range('D31').Value = '11.14'
$ws.Range('E31').Value = '  +1.22%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').Value = '556.83'
$ws.Range('E32').Value = '  -2.36%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').Value = '0.107'
$ws.Range('E33').Value = '  +1.58%  '
$ws.Range('B34').Value = 'dogwifhat'
$ws.Range('C34').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D34').Value = '3.61'
$ws.Range('E34').Value = '  +7.96%  '
$ws.Range('D35').Value = '58.61'
$ws.Range('E35').Value = '  +2.89%  '
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('D37').Value = '3.666.61'
$ws.Range('E37').Value = '  -2.45%  '
$ws.Range('D38').Value = '0.141'
$ws.Range('E38').Value = '  +4.78%  '
$ws.Range('D39').Value = '35.83'
$ws.Range('E39').Value = '  +0.91%  '
$ws.Range('D40').Value = '0.0₃0733'
$ws.Range('E40').Value = '  +8.30%  '
$ws.Range('D41').Value = '3.30'
$ws.Range('D42').Value = '2.70'
$ws.Range('E42').Value = '  +3.30%  '
$ws.Range('D43').Value = '0.0432'
$ws.Range('E43').Value = '  +5.88%  '
$ws.Range('D44').Value = '0.339'
$ws.Range('E44').Value = '  +1.66%  '
$ws.Range('D45').Value = '3.36'
$ws.Range('E45').Value = '  +1.80%  '
$ws.Range('D46').Value = '2.68'
$ws.Range('E46').Value = '  +1.64%  '
$ws.Range('E47').Value = '  +0.88%  '
$ws.Range('E48').Value = '  +4.88%  '
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('D50').Value = '129.88'
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('D51').Value = '2.66'
$ws.Range('E51').Value = '  -1.09%  '

# Restore default (General/Normal) formatting on the touched range.
$ws.Range("D2:E51").NumberFormat = "General"
$ws.Range("D2:E51").Style = "Normal"

